$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("H94").Value = 3250
$ws.Range("I94").Value = 3250
$ws.Range("K94").Value = 3250
$ws.Range("M94").Value = -2799
$ws.Range("H98").Value = 1998.75
$ws.Range("J98").Value = 1253
$ws.Range("L98").Value = 1253
$ws.Range("N98").Value = -4249
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 1998.75
$ws.Range("J122").Value = 1253
$ws.Range("L122").Value = 3759
$ws.Range("N122").Value = -8659

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1101.36
$ws.Range("I32").Value = 1111.4747
$ws.Range("J32").Value = 100
$ws.Range("K32").Value = 1111.4747
$ws.Range("L32").Value = 100
$ws.Range("M32").Value = -824.4747
$ws.Range("N32").Value = -674
$ws.Range("H61").Value = 1194.7241
$ws.Range("I61").Value = 1227.7291
$ws.Range("K61").Value = 1227.7291
$ws.Range("M61").Value = -1015.7291
$ws.Range("H63").Value = 2025.4333
$ws.Range("I63").Value = 1941.8695
$ws.Range("J63").Value = 2300
$ws.Range("K63").Value = 1941.8695
$ws.Range("L63").Value = 2300
$ws.Range("M63").Value = -1255.8695
$ws.Range("N63").Value = -3672
$ws.Range("H66").Value = 2025.4333
$ws.Range("I66").Value = 1941.8695
$ws.Range("J66").Value = 2300
$ws.Range("K66").Value = 9709.3475
$ws.Range("L66").Value = 11500
$ws.Range("M66").Value = -6277.3475
$ws.Range("N66").Value = -18364
$ws.Range("H74").Value = 925.04
$ws.Range("I74").Value = 903.3488
$ws.Range("J74").Value = 1058.2858
$ws.Range("K74").Value = 903.3488
$ws.Range("L74").Value = 1058.2858
$ws.Range("M74").Value = -29.34879999999998
$ws.Range("N74").Value = -2806.2858
$ws.Range("H77").Value = 925.04
$ws.Range("I77").Value = 903.3488
$ws.Range("J77").Value = 1058.2858
$ws.Range("K77").Value = 4516.744
$ws.Range("L77").Value = 5291.429
$ws.Range("M77").Value = -148.7439999999997
$ws.Range("N77").Value = -14027.429
$ws.Range("H132").Value = 1178.8833
$ws.Range("I132").Value = 965.9808
$ws.Range("J132").Value = 2562.75
$ws.Range("K132").Value = 2897.9424
$ws.Range("L132").Value = 7688.25
$ws.Range("M132").Value = -367.9423999999999
$ws.Range("N132").Value = -12748.25
$ws.Range("H136").Value = 1194.7241
$ws.Range("I136").Value = 1227.7291
$ws.Range("K136").Value = 3683.1873
$ws.Range("M136").Value = -1133.1873

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 10000
$ws.Range("I7").Value = 10000
$ws.Range("K7").Value = 10000
$ws.Range("M7").Value = -9887
$ws.Range("H134").Value = 32891.062
$ws.Range("I134").Value = 2068
$ws.Range("J134").Value = 171594.83
$ws.Range("K134").Value = 6204
$ws.Range("L134").Value = 514784.49
$ws.Range("M134").Value = -3669
$ws.Range("N134").Value = -519854.49

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 41608.348
$ws.Range("I31").Value = 3239.0557
$ws.Range("J31").Value = 127939.25
$ws.Range("K31").Value = 3239.0557
$ws.Range("L31").Value = 127939.25
$ws.Range("M31").Value = -2944.0557
$ws.Range("N31").Value = -128529.25
$ws.Range("H34").Value = 41608.348
$ws.Range("I34").Value = 3239.0557
$ws.Range("J34").Value = 127939.25
$ws.Range("K34").Value = 3239.0557
$ws.Range("L34").Value = 127939.25
$ws.Range("M34").Value = -3037.0557
$ws.Range("N34").Value = -128343.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 971.4167
$ws.Range("J121").Value = 971.4167
$ws.Range("L121").Value = 2914.2501
$ws.Range("N121").Value = -5534.2501
$ws.Range("H131").Value = 16701633
$ws.Range("J131").Value = 1527.56
$ws.Range("L131").Value = 4582.68
$ws.Range("N131").Value = -14662.68

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 39901.293
$ws.Range("J123").Value = 39901.293
$ws.Range("L123").Value = 39901.293
$ws.Range("N123").Value = -44801.293

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1416.3667
$ws.Range("I7").Value = 1388.2273
$ws.Range("K7").Value = 1388.2273
$ws.Range("M7").Value = -1276.2273
$ws.Range("H40").Value = 723184.9399999999
$ws.Range("I40").Value = 1011448.9
$ws.Range("J40").Value = 2525
$ws.Range("K40").Value = 1011448.9
$ws.Range("L40").Value = 2525
$ws.Range("M40").Value = -1011312.9
$ws.Range("N40").Value = -2797
$ws.Range("H55").Value = 262.1111
$ws.Range("I55").Value = 295.7
$ws.Range("K55").Value = 295.7
$ws.Range("M55").Value = -122.7
$ws.Range("H61").Value = 1963.2
$ws.Range("I61").Value = 2080.6155
$ws.Range("J61").Value = 1200
$ws.Range("K61").Value = 2080.6155
$ws.Range("L61").Value = 1200
$ws.Range("M61").Value = -1878.6155
$ws.Range("N61").Value = -1604
$ws.Range("H100").Value = 6946376
$ws.Range("I100").Value = 10102792
$ws.Range("J100").Value = 2260.8
$ws.Range("K100").Value = 10102792
$ws.Range("L100").Value = 2260.8
$ws.Range("M100").Value = -10102251
$ws.Range("N100").Value = -3342.8
$ws.Range("H101").Value = 20681
$ws.Range("J101").Value = 20681
$ws.Range("L101").Value = 20681
$ws.Range("N101").Value = -27171
$ws.Range("H113").Value = 1963.2
$ws.Range("I113").Value = 2080.6155
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 2080.6155
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 89.38450000000012
$ws.Range("N113").Value = -5540
$ws.Range("H126").Value = 1416.3667
$ws.Range("I126").Value = 1388.2273
$ws.Range("K126").Value = 4164.6819
$ws.Range("M126").Value = -1694.6819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 31663.334
$ws.Range("J94").Value = 31663.334
$ws.Range("L94").Value = 31663.334
$ws.Range("N94").Value = -33465.334
$ws.Range("H103").Value = 30602
$ws.Range("J103").Value = 30602
$ws.Range("L103").Value = 30602
$ws.Range("N103").Value = -32946
$ws.Range("H113").Value = 434.15384
$ws.Range("I113").Value = 433.25
$ws.Range("K113").Value = 1299.75
$ws.Range("M113").Value = 870.25
$ws.Range("H132").Value = 1132.4839
$ws.Range("I132").Value = 982.12
$ws.Range("J132").Value = 1759
$ws.Range("K132").Value = 2946.36
$ws.Range("L132").Value = 5277
$ws.Range("M132").Value = -416.3600000000001
$ws.Range("N132").Value = -10337
$ws.Range("H136").Value = 1309.8
$ws.Range("I136").Value = 1219.5428
$ws.Range("J136").Value = 1941.6
$ws.Range("K136").Value = 3658.6284
$ws.Range("L136").Value = 5824.799999999999
$ws.Range("M136").Value = -1108.6284
$ws.Range("N136").Value = -10924.8
